$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two brand-new rows (line7, line8) are added to the table; give their
# "id" cell (column A) the same look as the rest of the column (bold,
# centered, thin box border) by copying the format from an existing cell.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null

# Final data for the whole table (rows 2-17): id, name, from_bus, to_bus, in_service
$data = @(
    @(0,  "line1", 7,  9,  $true),
    @(1,  "line2", 9,  8,  $false),
    @(2,  "line3", 8,  10, $true),
    @(3,  "line4", 8,  11, $true),
    @(4,  "line5", 10, 5,  $true),
    @(5,  "line6", 12, 8,  $true),
    @(6,  "line7", 14, 11, $true),
    @(7,  "line8", 16, 9,  $true),
    @(8,  "extr1", 5,  12, $true),
    @(9,  "extr2", 5,  9,  $true),
    @(10, "extr3", 10, 11, $false),
    @(11, "extr4", 7,  8,  $true),
    @(12, "extr5", 9,  11, $true),
    @(13, "extr6", 7,  11, $false),
    @(14, "extr7", 5,  7,  $true),
    @(15, "extr8", 8,  5,  $false)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
